# Projects summary template: add a banner/title row and two new grouped
# sections (Year Specific Overview, O&M Costs Totals) to the header sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new row above the existing header row; the current header
#    (row 1) becomes row 2.
# ---------------------------------------------------------------------
$ws.Rows.Item(1).Insert()

# ---------------------------------------------------------------------
# 2. Rename a few header labels on (now) row 2 to their shorter forms.
# ---------------------------------------------------------------------
$ws.Range("C2").Value = "Title"
$ws.Range("D2").Value = "Overview"
$ws.Range("M2").Value = "Leads"

# ---------------------------------------------------------------------
# 3. Add the "Year Specific Overview" block of columns (R:AD), with a
#    narrow spacer column Q separating it from the original A:P block.
# ---------------------------------------------------------------------
$ws.Range("R2").Value = "Priorities"
$ws.Range("S2").Value = "Special Equipment?"
$ws.Range("T2").Value = "Field Component?"
$ws.Range("U2").Value = "Data Collected or Generated?"
$ws.Range("V2").Value = "Type of Data"
$ws.Range("W2").Value = "Data Products"
$ws.Range("X2").Value = "Data Management"
$ws.Range("Y2").Value = "Labortory Work?"
$ws.Range("Z2").Value = "Special IT Requirements"
$ws.Range("AA2").Value = "Additional Notes"
$ws.Range("AB2").Value = "Financial Coding"
$ws.Range("AC2").Value = "Submitted"
$ws.Range("AD2").Value = "Status"

# ---------------------------------------------------------------------
# 4. Add the "O&M Costs Totals" block of columns (AF:AI), with a narrow
#    spacer column AE separating it from the previous block.
# ---------------------------------------------------------------------
$ws.Range("AF2").Value = "A-Base"
$ws.Range("AG2").Value = "B-Base"
$ws.Range("AH2").Value = "C-Base"
$ws.Range("AI2").Value = "Total"

# ---------------------------------------------------------------------
# 5. Style the two narrow spacer columns (Q and AE) with a solid black
#    fill, and give them the template's narrow width.
# ---------------------------------------------------------------------
$ws.Range("Q1:Q2").Interior.ColorIndex = 1
$ws.Range("AE1:AE2").Interior.ColorIndex = 1
$ws.Columns.Item(17).ColumnWidth = 1.75
$ws.Columns.Item(31).ColumnWidth = 1.75

# ---------------------------------------------------------------------
# 6. Give the new data columns their "best fit" widths, matching the
#    template's general look (auto-sized to content).
# ---------------------------------------------------------------------
$ws.Columns.Item(18).ColumnWidth = 8.92
$ws.Columns.Item(19).ColumnWidth = 19.59
$ws.Columns.Item(20).ColumnWidth = 17.92
$ws.Columns.Item(21).ColumnWidth = 29.59
$ws.Columns.Item(22).ColumnWidth = 12.59
$ws.Columns.Item(23).ColumnWidth = 14.09
$ws.Columns.Item(24).ColumnWidth = 18.59
$ws.Columns.Item(25).ColumnWidth = 17.09
$ws.Columns.Item(26).ColumnWidth = 23.75
$ws.Columns.Item(27).ColumnWidth = 16.92
$ws.Columns.Item(28).ColumnWidth = 16.09
$ws.Columns.Item(29).ColumnWidth = 10.42
$ws.Columns.Item(30).ColumnWidth = 6.26

# ---------------------------------------------------------------------
# 7. Build the big banner row (row 1): three merged, dark-filled, bold
#    white-centered title cells spanning each block of columns.
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 30.75

$ws.Range("A1").Value = "Project General Overview"
$ws.Range("R1").Value = "Year Specific Overview"
$ws.Range("AF1").Value = "O&M Costs Totals"

# Build the title style once on A1, then copy its *formatting only* onto
# the other two banner cells (keeps the generated style table small).
$titleCell = $ws.Range("A1")
$titleCell.Style = "Accent5"
$titleCell.Font.Bold = $true
$titleCell.Font.Size = 16
$titleCell.HorizontalAlignment = -4108

$titleCell.Copy()
$ws.Range("R1").PasteSpecial(-4122)
$ws.Range("AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 8. Merge each banner title across its block of columns.
# ---------------------------------------------------------------------
$ws.Range("A1:P1").Merge()
$ws.Range("R1:AD1").Merge()
$ws.Range("AF1:AI1").Merge()

# ---------------------------------------------------------------------
# 9. Selection / dimension bookkeeping to match the finished template.
# ---------------------------------------------------------------------
$ws.Range("A1:P1").Select()
